$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "3A-DSM"
$ws.Range("B8").Value = "2-Jun"
$ws.Range("C8").Value = "24-Jun"
$ws.Range("D8").Value = "25-Jun"
$ws.Range("E8").Value = "9-Jul"
$ws.Range("F8").Value = "23-Jul"
$ws.Range("G8").Value = "24-Jul"
$ws.Range("A9").Value = "PROBABILIDAD Y ESTADÍSTICA"
$ws.Range("B9").Value = "Martes 26 de Mayo `n Aula C-20 - 7:00"
$ws.Range("C9").Value = "Martes 23 de Junio `n Aula C-20 - 7:00"
$ws.Range("D9").Value = "Martes 23 de Junio `n Aula C-20 - 7:00"
$ws.Range("E9").Value = "Martes 7 de Julio `n Aula C-20 - 7:00"
$ws.Range("F9").Value = "Miercoles 22 de Julio `n Lab4 - 11:00"
$ws.Range("G9").Value = "Jueves 23 de Julio `n Aula C-17 - 12:40"
$ws.Range("A10").Value = "FORMACIÓN SOCIOCULTURAL III"
$ws.Range("B10").Value = "Jueves 28 de Mayo `n Aula C-16 - 7:00"
$ws.Range("C10").Value = "Jueves 18 de Junio `n Aula C-16 - 7:00"
$ws.Range("D10").Value = "Jueves 18 de Junio `n Aula C-16 - 7:00"
$ws.Range("E10").Value = "Jueves 2 de Julio `n Aula C-16 - 7:00"
$ws.Range("F10").Value = "Miercoles 22 de Julio `n Aula C-23 - 7:00"
$ws.Range("G10").Value = "Jueves 23 de Julio `n Aula C-16 - 7:00"
$ws.Range("A11").Value = "CALCULO DIFERENCIAL"
$ws.Range("B11").Value = "Viernes 29 de Mayo `n Aula C-17 - 7:00"
$ws.Range("C11").Value = "Viernes 19 de Junio `n Aula C-17 - 7:00"
$ws.Range("D11").Value = "Viernes 19 de Junio `n Aula C-17 - 7:00"
$ws.Range("E11").Value = "Viernes 3 de Julio `n Aula C-17 - 7:00"
$ws.Range("F11").Value = "Miercoles 22 de Julio `n Lab8 - 9:00"
$ws.Range("G11").Value = "Jueves 23 de Julio `n Lab6 - 9:00"
$ws.Range("A12").Value = "APLICACIONES WEB"
$ws.Range("B12").Value = "Miercoles 27 de Mayo `n Lab4 - 11:00"
$ws.Range("C12").Value = "Miercoles 17 de Junio `n Lab4 - 11:00"
$ws.Range("D12").Value = "Miercoles 24 de Junio `n Lab4 - 11:00"
$ws.Range("E12").Value = "Miercoles 8 de Julio `n Lab4 - 11:00"
$ws.Range("F12").Value = "Martes 21 de Julio `n Aula C-20 - 7:00"
$ws.Range("G12").Value = "Miercoles 22 de Julio `n Aula C-23 - 7:00"
$ws.Range("A13").Value = "INTEGRADORA I"
$ws.Range("B13").Value = "Viernes 29 de Mayo `n Lab3 - 11:00"
$ws.Range("C13").Value = "Viernes 19 de Junio `n Lab3 - 11:00"
$ws.Range("D13").Value = "Viernes 19 de Junio `n Lab3 - 11:00"
$ws.Range("E13").Value = "Viernes 3 de Julio `n Lab3 - 11:00"
$ws.Range("F13").Value = "Miercoles 22 de Julio `n Aula C-23 - 7:50"
$ws.Range("G13").Value = "Jueves 23 de Julio `n Aula C-16 - 7:50"
$ws.Range("A14").Value = "SISTEMAS OPERATIVOS"
$ws.Range("B14").Value = "Lunes 1 de Junio `n Lab8 - 9:00"
$ws.Range("C14").Value = "Miercoles 17 de Junio `n Lab8 - 9:00"
$ws.Range("D14").Value = "Lunes 22 de Junio `n Lab8 - 9:00"
$ws.Range("E14").Value = "Lunes 6 de Julio `n Lab8 - 9:00"
$ws.Range("F14").Value = "Miercoles 22 de Julio `n Lab8 - 9:50"
$ws.Range("G14").Value = "Jueves 23 de Julio `n Lab6 - 9:50"
$ws.Range("A15").Value = "BASES DE DATOS PARA APLICACIONES"
$ws.Range("B15").Value = "Lunes 1 de Junio `n Lab2 - 11:00"
$ws.Range("C15").Value = "Lunes 22 de Junio `n Lab2 - 11:00"
$ws.Range("D15").Value = "Lunes 22 de Junio `n Lab2 - 11:00"
$ws.Range("E15").Value = "Lunes 6 de Julio `n Lab2 - 11:00"
$ws.Range("F15").Value = "Miercoles 22 de Julio `n Lab4 - 11:50"
$ws.Range("G15").Value = "Jueves 23 de Julio `n Lab2 - 13:30"
$ws.Range("A16").Value = "INGLÉS III"
